$wb = $excel.ActiveWorkbook

# ---- Sheet: Crystal Arcade ----
$ws = $wb.Worksheets.Item('Crystal Arcade')
$ws.Range('A64:N64').Copy($ws.Range('A65:N65'))
$ws.Cells.Item(65, 1).Value2 = 'LUMI'
$ws.Cells.Item(65, 2).Value2 = 'FRANK'
$ws.Cells.Item(65, 3).Value2 = 'CHARLIE'
$ws.Cells.Item(65, 4).Value2 = 'TARA'
$ws.Cells.Item(65, 5).Value2 = 'DRACO'
$ws.Cells.Item(65, 6).Value2 = 'GRAY'
$ws.Cells.Item(65, 7).Value2 = 'Equipo 2'
$ws.Cells.Item(65, 8).Value2 = 'ER|Naipishu😎'
$ws.Cells.Item(65, 9).Value2 = 'あの頃のしずく👍'
$ws.Cells.Item(65, 10).Value2 = 'ZETA|Levi'
$ws.Cells.Item(65, 11).Value2 = 'RVL|Terry'
$ws.Cells.Item(65, 12).Value2 = 'RVL|I see'
$ws.Cells.Item(65, 13).Value2 = 'Mameshi'
$ws.Cells.Item(65, 14).Value2 = '20250726T133554.000Z'

$ws.Range('A65:N65').Copy($ws.Range('A66:N66'))
$ws.Range('G63').Copy($ws.Range('G66'))
$ws.Cells.Item(66, 1).Value2 = 'SANDY'
$ws.Cells.Item(66, 2).Value2 = 'LILY'
$ws.Cells.Item(66, 3).Value2 = 'GUS'
$ws.Cells.Item(66, 4).Value2 = 'MOE'
$ws.Cells.Item(66, 5).Value2 = 'CHARLIE'
$ws.Cells.Item(66, 6).Value2 = 'BULL'
$ws.Cells.Item(66, 7).Value2 = 'Equipo 1'
$ws.Cells.Item(66, 8).Value2 = 'ER|Naipishu😎'
$ws.Cells.Item(66, 9).Value2 = 'あの頃のしずく👍'
$ws.Cells.Item(66, 10).Value2 = 'ZETA|Levi'
$ws.Cells.Item(66, 11).Value2 = 'RVL|Terry'
$ws.Cells.Item(66, 12).Value2 = 'Mameshi'
$ws.Cells.Item(66, 13).Value2 = 'RVL|I see'
$ws.Cells.Item(66, 14).Value2 = '20250726T132900.000Z'

$ws.Range('A66:N66').Copy($ws.Range('A67:N67'))
$ws.Cells.Item(67, 1).Value2 = 'SANDY'
$ws.Cells.Item(67, 2).Value2 = 'LILY'
$ws.Cells.Item(67, 3).Value2 = 'GUS'
$ws.Cells.Item(67, 4).Value2 = 'MOE'
$ws.Cells.Item(67, 5).Value2 = 'CHARLIE'
$ws.Cells.Item(67, 6).Value2 = 'BULL'
$ws.Cells.Item(67, 7).Value2 = 'Equipo 1'
$ws.Cells.Item(67, 8).Value2 = 'ER|Naipishu😎'
$ws.Cells.Item(67, 9).Value2 = 'あの頃のしずく👍'
$ws.Cells.Item(67, 10).Value2 = 'ZETA|Levi'
$ws.Cells.Item(67, 11).Value2 = 'RVL|Terry'
$ws.Cells.Item(67, 12).Value2 = 'Mameshi'
$ws.Cells.Item(67, 13).Value2 = 'RVL|I see'
$ws.Cells.Item(67, 14).Value2 = '20250726T132657.000Z'

# ---- Sheet: Hot Potato ----
$ws = $wb.Worksheets.Item('Hot Potato')
$ws.Range('A80:N80').Copy($ws.Range('A81:N81'))
$ws.Cells.Item(81, 1).Value2 = 'EDGAR'
$ws.Cells.Item(81, 2).Value2 = 'KIT'
$ws.Cells.Item(81, 3).Value2 = 'BONNIE'
$ws.Cells.Item(81, 4).Value2 = 'FANG'
$ws.Cells.Item(81, 5).Value2 = 'DRACO'
$ws.Cells.Item(81, 6).Value2 = 'CORDELIUS'
$ws.Cells.Item(81, 7).Value2 = 'Equipo 2'
$ws.Cells.Item(81, 8).Value2 = 'ER|Naipishu😎'
$ws.Cells.Item(81, 9).Value2 = 'あの頃のしずく👍'
$ws.Cells.Item(81, 10).Value2 = 'ZETA|Levi'
$ws.Cells.Item(81, 11).Value2 = 'RVL|Terry'
$ws.Cells.Item(81, 12).Value2 = 'RVL|I see'
$ws.Cells.Item(81, 13).Value2 = 'Mameshi'
$ws.Cells.Item(81, 14).Value2 = '20250726T132043.000Z'

$ws.Range('A81:N81').Copy($ws.Range('A82:N82'))
$ws.Cells.Item(82, 1).Value2 = 'EDGAR'
$ws.Cells.Item(82, 2).Value2 = 'KIT'
$ws.Cells.Item(82, 3).Value2 = 'BONNIE'
$ws.Cells.Item(82, 4).Value2 = 'FANG'
$ws.Cells.Item(82, 5).Value2 = 'DRACO'
$ws.Cells.Item(82, 6).Value2 = 'CORDELIUS'
$ws.Cells.Item(82, 7).Value2 = 'Equipo 2'
$ws.Cells.Item(82, 8).Value2 = 'ER|Naipishu😎'
$ws.Cells.Item(82, 9).Value2 = 'あの頃のしずく👍'
$ws.Cells.Item(82, 10).Value2 = 'ZETA|Levi'
$ws.Cells.Item(82, 11).Value2 = 'RVL|Terry'
$ws.Cells.Item(82, 12).Value2 = 'RVL|I see'
$ws.Cells.Item(82, 13).Value2 = 'Mameshi'
$ws.Cells.Item(82, 14).Value2 = '20250726T131852.000Z'

$ws.Range('A82:N82').Copy($ws.Range('A83:N83'))
$ws.Cells.Item(83, 1).Value2 = 'LILY'
$ws.Cells.Item(83, 2).Value2 = 'BUZZ'
$ws.Cells.Item(83, 3).Value2 = 'CHARLIE'
$ws.Cells.Item(83, 4).Value2 = 'AMBER'
$ws.Cells.Item(83, 5).Value2 = 'CARL'
$ws.Cells.Item(83, 6).Value2 = 'KIT'
$ws.Cells.Item(83, 7).Value2 = 'Equipo 2'
$ws.Cells.Item(83, 8).Value2 = 'ER|Naipishu😎'
$ws.Cells.Item(83, 9).Value2 = 'あの頃のしずく👍'
$ws.Cells.Item(83, 10).Value2 = 'ZETA|Levi'
$ws.Cells.Item(83, 11).Value2 = 'RVL|Terry'
$ws.Cells.Item(83, 12).Value2 = 'RVL|I see'
$ws.Cells.Item(83, 13).Value2 = 'Mameshi'
$ws.Cells.Item(83, 14).Value2 = '20250726T131220.000Z'

$ws.Range('A83:N83').Copy($ws.Range('A84:N84'))
$ws.Cells.Item(84, 1).Value2 = 'LILY'
$ws.Cells.Item(84, 2).Value2 = 'BUZZ'
$ws.Cells.Item(84, 3).Value2 = 'CHARLIE'
$ws.Cells.Item(84, 4).Value2 = 'AMBER'
$ws.Cells.Item(84, 5).Value2 = 'CARL'
$ws.Cells.Item(84, 6).Value2 = 'KIT'
$ws.Cells.Item(84, 7).Value2 = 'Equipo 2'
$ws.Cells.Item(84, 8).Value2 = 'ER|Naipishu😎'
$ws.Cells.Item(84, 9).Value2 = 'あの頃のしずく👍'
$ws.Cells.Item(84, 10).Value2 = 'ZETA|Levi'
$ws.Cells.Item(84, 11).Value2 = 'RVL|Terry'
$ws.Cells.Item(84, 12).Value2 = 'RVL|I see'
$ws.Cells.Item(84, 13).Value2 = 'Mameshi'
$ws.Cells.Item(84, 14).Value2 = '20250726T131015.000Z'

# ---- Sheet: Layer Cake ----
$ws = $wb.Worksheets.Item('Layer Cake')
$ws.Range('A85:N85').Copy($ws.Range('A86:N86'))
$ws.Cells.Item(86, 1).Value2 = 'CORDELIUS'
$ws.Cells.Item(86, 2).Value2 = 'KIT'
$ws.Cells.Item(86, 3).Value2 = 'DRACO'
$ws.Cells.Item(86, 4).Value2 = 'BUSTER'
$ws.Cells.Item(86, 5).Value2 = 'CHARLIE'
$ws.Cells.Item(86, 6).Value2 = 'JAE-YONG'
$ws.Cells.Item(86, 7).Value2 = 'Equipo 1'
$ws.Cells.Item(86, 8).Value2 = 'CR|Milkreo'
$ws.Cells.Item(86, 9).Value2 = 'CR|Tensai'
$ws.Cells.Item(86, 10).Value2 = 'CR|Moya'
$ws.Cells.Item(86, 11).Value2 = 'NAVI|Ryohei'
$ws.Cells.Item(86, 12).Value2 = 'NAVI|Achapi'
$ws.Cells.Item(86, 13).Value2 = 'NAVI|Kuru'
$ws.Cells.Item(86, 14).Value2 = '20250726T133545.000Z'

$ws.Range('A86:N86').Copy($ws.Range('A87:N87'))
$ws.Cells.Item(87, 1).Value2 = 'CORDELIUS'
$ws.Cells.Item(87, 2).Value2 = 'KIT'
$ws.Cells.Item(87, 3).Value2 = 'DRACO'
$ws.Cells.Item(87, 4).Value2 = 'BUSTER'
$ws.Cells.Item(87, 5).Value2 = 'CHARLIE'
$ws.Cells.Item(87, 6).Value2 = 'JAE-YONG'
$ws.Cells.Item(87, 7).Value2 = 'Equipo 1'
$ws.Cells.Item(87, 8).Value2 = 'CR|Milkreo'
$ws.Cells.Item(87, 9).Value2 = 'CR|Tensai'
$ws.Cells.Item(87, 10).Value2 = 'CR|Moya'
$ws.Cells.Item(87, 11).Value2 = 'NAVI|Ryohei'
$ws.Cells.Item(87, 12).Value2 = 'NAVI|Achapi'
$ws.Cells.Item(87, 13).Value2 = 'NAVI|Kuru'
$ws.Cells.Item(87, 14).Value2 = '20250726T133325.000Z'

$ws.Range('A87:N87').Copy($ws.Range('A88:N88'))
$ws.Range('G84').Copy($ws.Range('G88'))
$ws.Cells.Item(88, 1).Value2 = 'CORDELIUS'
$ws.Cells.Item(88, 2).Value2 = 'KIT'
$ws.Cells.Item(88, 3).Value2 = 'DRACO'
$ws.Cells.Item(88, 4).Value2 = 'BUSTER'
$ws.Cells.Item(88, 5).Value2 = 'CHARLIE'
$ws.Cells.Item(88, 6).Value2 = 'JAE-YONG'
$ws.Cells.Item(88, 7).Value2 = 'Equipo 2'
$ws.Cells.Item(88, 8).Value2 = 'CR|Milkreo'
$ws.Cells.Item(88, 9).Value2 = 'CR|Tensai'
$ws.Cells.Item(88, 10).Value2 = 'CR|Moya'
$ws.Cells.Item(88, 11).Value2 = 'NAVI|Ryohei'
$ws.Cells.Item(88, 12).Value2 = 'NAVI|Achapi'
$ws.Cells.Item(88, 13).Value2 = 'NAVI|Kuru'
$ws.Cells.Item(88, 14).Value2 = '20250726T133106.000Z'

# ---- Sheet: Open Business ----
$ws = $wb.Worksheets.Item('Open Business')
$ws.Range('A97:N97').Copy($ws.Range('A98:N98'))
$ws.Range('G96').Copy($ws.Range('G98'))
$ws.Cells.Item(98, 1).Value2 = 'SANDY'
$ws.Cells.Item(98, 2).Value2 = 'KAZE'
$ws.Cells.Item(98, 3).Value2 = 'CROW'
$ws.Cells.Item(98, 4).Value2 = 'SHADE'
$ws.Cells.Item(98, 5).Value2 = 'CORDELIUS'
$ws.Cells.Item(98, 6).Value2 = 'AMBER'
$ws.Cells.Item(98, 7).Value2 = 'Equipo 2'
$ws.Cells.Item(98, 8).Value2 = 'CR|Milkreo'
$ws.Cells.Item(98, 9).Value2 = 'CR|Tensai'
$ws.Cells.Item(98, 10).Value2 = 'CR|Moya'
$ws.Cells.Item(98, 11).Value2 = 'NAVI|Achapi'
$ws.Cells.Item(98, 12).Value2 = 'NAVI|Ryohei'
$ws.Cells.Item(98, 13).Value2 = 'NAVI|Kuru'
$ws.Cells.Item(98, 14).Value2 = '20250726T132427.000Z'

$ws.Range('A98:N98').Copy($ws.Range('A99:N99'))
$ws.Cells.Item(99, 1).Value2 = 'SANDY'
$ws.Cells.Item(99, 2).Value2 = 'KAZE'
$ws.Cells.Item(99, 3).Value2 = 'CROW'
$ws.Cells.Item(99, 4).Value2 = 'SHADE'
$ws.Cells.Item(99, 5).Value2 = 'CORDELIUS'
$ws.Cells.Item(99, 6).Value2 = 'AMBER'
$ws.Cells.Item(99, 7).Value2 = 'Equipo 2'
$ws.Cells.Item(99, 8).Value2 = 'CR|Milkreo'
$ws.Cells.Item(99, 9).Value2 = 'CR|Tensai'
$ws.Cells.Item(99, 10).Value2 = 'CR|Moya'
$ws.Cells.Item(99, 11).Value2 = 'NAVI|Achapi'
$ws.Cells.Item(99, 12).Value2 = 'NAVI|Ryohei'
$ws.Cells.Item(99, 13).Value2 = 'NAVI|Kuru'
$ws.Cells.Item(99, 14).Value2 = '20250726T132157.000Z'

$ws.Range('A99:N99').Copy($ws.Range('A100:N100'))
$ws.Cells.Item(100, 1).Value2 = 'EMZ'
$ws.Cells.Item(100, 2).Value2 = 'KAZE'
$ws.Cells.Item(100, 3).Value2 = 'ASH'
$ws.Cells.Item(100, 4).Value2 = 'BUZZ'
$ws.Cells.Item(100, 5).Value2 = 'SANDY'
$ws.Cells.Item(100, 6).Value2 = 'MEG'
$ws.Cells.Item(100, 7).Value2 = 'Equipo 2'
$ws.Cells.Item(100, 8).Value2 = 'CR|Milkreo'
$ws.Cells.Item(100, 9).Value2 = 'CR|Tensai'
$ws.Cells.Item(100, 10).Value2 = 'CR|Moya'
$ws.Cells.Item(100, 11).Value2 = 'NAVI|Ryohei'
$ws.Cells.Item(100, 12).Value2 = 'NAVI|Achapi'
$ws.Cells.Item(100, 13).Value2 = 'NAVI|Kuru'
$ws.Cells.Item(100, 14).Value2 = '20250726T131616.000Z'

$ws.Range('A100:N100').Copy($ws.Range('A101:N101'))
$ws.Cells.Item(101, 1).Value2 = 'EMZ'
$ws.Cells.Item(101, 2).Value2 = 'KAZE'
$ws.Cells.Item(101, 3).Value2 = 'ASH'
$ws.Cells.Item(101, 4).Value2 = 'BUZZ'
$ws.Cells.Item(101, 5).Value2 = 'SANDY'
$ws.Cells.Item(101, 6).Value2 = 'MEG'
$ws.Cells.Item(101, 7).Value2 = 'Equipo 2'
$ws.Cells.Item(101, 8).Value2 = 'CR|Milkreo'
$ws.Cells.Item(101, 9).Value2 = 'CR|Tensai'
$ws.Cells.Item(101, 10).Value2 = 'CR|Moya'
$ws.Cells.Item(101, 11).Value2 = 'NAVI|Ryohei'
$ws.Cells.Item(101, 12).Value2 = 'NAVI|Achapi'
$ws.Cells.Item(101, 13).Value2 = 'NAVI|Kuru'
$ws.Cells.Item(101, 14).Value2 = '20250726T131428.000Z'

# ---- Sheet: Dry Season ----
$ws = $wb.Worksheets.Item('Dry Season')
$ws.Range('A60:N60').Copy($ws.Range('A61:N61'))
$ws.Cells.Item(61, 1).Value2 = 'KIT'
$ws.Cells.Item(61, 2).Value2 = 'R-T'
$ws.Cells.Item(61, 3).Value2 = 'SPIKE'
$ws.Cells.Item(61, 4).Value2 = 'DOUG'
$ws.Cells.Item(61, 5).Value2 = 'BROCK'
$ws.Cells.Item(61, 6).Value2 = '8-BIT'
$ws.Cells.Item(61, 7).Value2 = 'Equipo 2'
$ws.Cells.Item(61, 8).Value2 = 'ER|Wahochi'
$ws.Cells.Item(61, 9).Value2 = 'FG|Rujao'
$ws.Cells.Item(61, 10).Value2 = 'FG|Nem🌙·̩͙⋆͛'
$ws.Cells.Item(61, 11).Value2 = 'FZ|Toridesu'
$ws.Cells.Item(61, 12).Value2 = 'FZ|Danshari'
$ws.Cells.Item(61, 13).Value2 = 'FZ|Mira'
$ws.Cells.Item(61, 14).Value2 = '20250726T133300.000Z'

$ws.Range('A61:N61').Copy($ws.Range('A62:N62'))
$ws.Cells.Item(62, 1).Value2 = 'KIT'
$ws.Cells.Item(62, 2).Value2 = 'R-T'
$ws.Cells.Item(62, 3).Value2 = 'SPIKE'
$ws.Cells.Item(62, 4).Value2 = 'DOUG'
$ws.Cells.Item(62, 5).Value2 = 'BROCK'
$ws.Cells.Item(62, 6).Value2 = '8-BIT'
$ws.Cells.Item(62, 7).Value2 = 'Equipo 2'
$ws.Cells.Item(62, 8).Value2 = 'ER|Wahochi'
$ws.Cells.Item(62, 9).Value2 = 'FG|Rujao'
$ws.Cells.Item(62, 10).Value2 = 'FG|Nem🌙·̩͙⋆͛'
$ws.Cells.Item(62, 11).Value2 = 'FZ|Toridesu'
$ws.Cells.Item(62, 12).Value2 = 'FZ|Danshari'
$ws.Cells.Item(62, 13).Value2 = 'FZ|Mira'
$ws.Cells.Item(62, 14).Value2 = '20250726T133040.000Z'

$ws.Range('A62:N62').Copy($ws.Range('A63:N63'))
$ws.Range('G59').Copy($ws.Range('G63'))
$ws.Cells.Item(63, 1).Value2 = 'KIT'
$ws.Cells.Item(63, 2).Value2 = 'R-T'
$ws.Cells.Item(63, 3).Value2 = 'SPIKE'
$ws.Cells.Item(63, 4).Value2 = 'DOUG'
$ws.Cells.Item(63, 5).Value2 = 'BROCK'
$ws.Cells.Item(63, 6).Value2 = '8-BIT'
$ws.Cells.Item(63, 7).Value2 = 'Equipo 1'
$ws.Cells.Item(63, 8).Value2 = 'ER|Wahochi'
$ws.Cells.Item(63, 9).Value2 = 'FG|Rujao'
$ws.Cells.Item(63, 10).Value2 = 'FG|Nem🌙·̩͙⋆͛'
$ws.Cells.Item(63, 11).Value2 = 'FZ|Toridesu'
$ws.Cells.Item(63, 12).Value2 = 'FZ|Danshari'
$ws.Cells.Item(63, 13).Value2 = 'FZ|Mira'
$ws.Cells.Item(63, 14).Value2 = '20250726T132821.000Z'

# ---- Sheet: Ring of Fire ----
$ws = $wb.Worksheets.Item('Ring of Fire')
$ws.Range('A81:N81').Copy($ws.Range('A82:N82'))
$ws.Cells.Item(82, 1).Value2 = 'BERRY'
$ws.Cells.Item(82, 2).Value2 = 'HANK'
$ws.Cells.Item(82, 3).Value2 = 'CROW'
$ws.Cells.Item(82, 4).Value2 = 'BEA'
$ws.Cells.Item(82, 5).Value2 = 'DOUG'
$ws.Cells.Item(82, 6).Value2 = 'PAM'
$ws.Cells.Item(82, 7).Value2 = 'Equipo 1'
$ws.Cells.Item(82, 8).Value2 = 'FG|Rujao'
$ws.Cells.Item(82, 9).Value2 = 'ER|Wahochi'
$ws.Cells.Item(82, 10).Value2 = 'FG|Nem🌙·̩͙⋆͛'
$ws.Cells.Item(82, 11).Value2 = 'FZ|Mira'
$ws.Cells.Item(82, 12).Value2 = 'FZ|Toridesu'
$ws.Cells.Item(82, 13).Value2 = 'FZ|Danshari'
$ws.Cells.Item(82, 14).Value2 = '20250726T132134.000Z'

$ws.Range('A82:N82').Copy($ws.Range('A83:N83'))
$ws.Range('G79').Copy($ws.Range('G83'))
$ws.Cells.Item(83, 1).Value2 = 'BERRY'
$ws.Cells.Item(83, 2).Value2 = 'HANK'
$ws.Cells.Item(83, 3).Value2 = 'CROW'
$ws.Cells.Item(83, 4).Value2 = 'BEA'
$ws.Cells.Item(83, 5).Value2 = 'DOUG'
$ws.Cells.Item(83, 6).Value2 = 'PAM'
$ws.Cells.Item(83, 7).Value2 = 'Equipo 2'
$ws.Cells.Item(83, 8).Value2 = 'FG|Rujao'
$ws.Cells.Item(83, 9).Value2 = 'ER|Wahochi'
$ws.Cells.Item(83, 10).Value2 = 'FG|Nem🌙·̩͙⋆͛'
$ws.Cells.Item(83, 11).Value2 = 'FZ|Mira'
$ws.Cells.Item(83, 12).Value2 = 'FZ|Toridesu'
$ws.Cells.Item(83, 13).Value2 = 'FZ|Danshari'
$ws.Cells.Item(83, 14).Value2 = '20250726T131935.000Z'

$ws.Range('A83:N83').Copy($ws.Range('A84:N84'))
$ws.Range('G82').Copy($ws.Range('G84'))
$ws.Cells.Item(84, 1).Value2 = 'BERRY'
$ws.Cells.Item(84, 2).Value2 = 'HANK'
$ws.Cells.Item(84, 3).Value2 = 'CROW'
$ws.Cells.Item(84, 4).Value2 = 'BEA'
$ws.Cells.Item(84, 5).Value2 = 'DOUG'
$ws.Cells.Item(84, 6).Value2 = 'PAM'
$ws.Cells.Item(84, 7).Value2 = 'Equipo 1'
$ws.Cells.Item(84, 8).Value2 = 'FG|Rujao'
$ws.Cells.Item(84, 9).Value2 = 'ER|Wahochi'
$ws.Cells.Item(84, 10).Value2 = 'FG|Nem🌙·̩͙⋆͛'
$ws.Cells.Item(84, 11).Value2 = 'FZ|Mira'
$ws.Cells.Item(84, 12).Value2 = 'FZ|Toridesu'
$ws.Cells.Item(84, 13).Value2 = 'FZ|Danshari'
$ws.Cells.Item(84, 14).Value2 = '20250726T131722.000Z'

$ws.Range('A84:N84').Copy($ws.Range('A85:N85'))
$ws.Range('G83').Copy($ws.Range('G85'))
$ws.Cells.Item(85, 1).Value2 = 'PENNY'
$ws.Cells.Item(85, 2).Value2 = 'MEG'
$ws.Cells.Item(85, 3).Value2 = 'AMBER'
$ws.Cells.Item(85, 4).Value2 = 'PAM'
$ws.Cells.Item(85, 5).Value2 = 'KAZE'
$ws.Cells.Item(85, 6).Value2 = 'CHARLIE'
$ws.Cells.Item(85, 7).Value2 = 'Equipo 2'
$ws.Cells.Item(85, 8).Value2 = 'ER|Wahochi'
$ws.Cells.Item(85, 9).Value2 = 'FG|Nem🌙·̩͙⋆͛'
$ws.Cells.Item(85, 10).Value2 = 'FG|Rujao'
$ws.Cells.Item(85, 11).Value2 = 'FZ|Mira'
$ws.Cells.Item(85, 12).Value2 = 'FZ|Toridesu'
$ws.Cells.Item(85, 13).Value2 = 'FZ|Danshari'
$ws.Cells.Item(85, 14).Value2 = '20250726T131048.000Z'

$ws.Range('A85:N85').Copy($ws.Range('A86:N86'))
$ws.Range('G84').Copy($ws.Range('G86'))
$ws.Cells.Item(86, 1).Value2 = 'PENNY'
$ws.Cells.Item(86, 2).Value2 = 'MEG'
$ws.Cells.Item(86, 3).Value2 = 'AMBER'
$ws.Cells.Item(86, 4).Value2 = 'PAM'
$ws.Cells.Item(86, 5).Value2 = 'KAZE'
$ws.Cells.Item(86, 6).Value2 = 'CHARLIE'
$ws.Cells.Item(86, 7).Value2 = 'Equipo 1'
$ws.Cells.Item(86, 8).Value2 = 'ER|Wahochi'
$ws.Cells.Item(86, 9).Value2 = 'FG|Nem🌙·̩͙⋆͛'
$ws.Cells.Item(86, 10).Value2 = 'FG|Rujao'
$ws.Cells.Item(86, 11).Value2 = 'FZ|Mira'
$ws.Cells.Item(86, 12).Value2 = 'FZ|Toridesu'
$ws.Cells.Item(86, 13).Value2 = 'FZ|Danshari'
$ws.Cells.Item(86, 14).Value2 = '20250726T130820.000Z'

$ws.Range('A86:N86').Copy($ws.Range('A87:N87'))
$ws.Range('G85').Copy($ws.Range('G87'))
$ws.Cells.Item(87, 1).Value2 = 'PENNY'
$ws.Cells.Item(87, 2).Value2 = 'MEG'
$ws.Cells.Item(87, 3).Value2 = 'AMBER'
$ws.Cells.Item(87, 4).Value2 = 'PAM'
$ws.Cells.Item(87, 5).Value2 = 'KAZE'
$ws.Cells.Item(87, 6).Value2 = 'CHARLIE'
$ws.Cells.Item(87, 7).Value2 = 'Equipo 2'
$ws.Cells.Item(87, 8).Value2 = 'ER|Wahochi'
$ws.Cells.Item(87, 9).Value2 = 'FG|Nem🌙·̩͙⋆͛'
$ws.Cells.Item(87, 10).Value2 = 'FG|Rujao'
$ws.Cells.Item(87, 11).Value2 = 'FZ|Mira'
$ws.Cells.Item(87, 12).Value2 = 'FZ|Toridesu'
$ws.Cells.Item(87, 13).Value2 = 'FZ|Danshari'
$ws.Cells.Item(87, 14).Value2 = '20250726T130615.000Z'

Write-Host 'Done applying scrims update'
